# Apply the "Updated analysis, re-structured files" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in place ---

# Row 2: SouthKorea - refreshed CFR2 / DiffDE values
$ws.Range("D2").Value = 0.02225547035720965
$ws.Range("E2").Value = -0.02225547035720965

# Row 5: USA (All) - refreshed CFR2 / DiffDE values
$ws.Range("D5").Value = 0.0559832330827703
$ws.Range("E5").Value = -0.0559832330827703

# --- Insert a new row for "USA / NYC" right after the USA/All row (row 5), ---
# --- pushing Spain (was row 6) and Italy (was row 7) down by one row.     ---
$ws.Range("A6").EntireRow.Insert()

$ws.Range("A6").Value = "USA"
$ws.Range("B6").Value = "NYC"
$ws.Range("C6").Value = 43943
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D6").Value = 0.07258978752642781
$ws.Range("E6").Value = -0.07258978752642781
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# --- Spain now lives on row 7 - refresh its CFR2 / DiffDE values ---
$ws.Range("D7").Value = 0.1040126934054379
$ws.Range("E7").Value = -0.1040126934054379

# --- Italy now lives on row 8 - refresh its CFR2 / DiffDE values ---
$ws.Range("D8").Value = 0.1300061819994451
$ws.Range("E8").Value = -0.1300061819994451
